$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 170506
$ws.Range("C4").Value = 161329
$ws.Range("C7").Value = 5.38
$ws.Range("C8").Value = 65.66
